$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had a "D2" column (header in K1) with no matching
# "D1"/"D3" columns. Insert a new column at K so the existing "D2" header +
# data slide over to column L, leaving room for a new "D1" column at K and
# appending a new "D3" column at M.
$ws.Range("K1").EntireColumn.Insert()

# The inserted column only carries a header value (row 1) - rows 2-4 have no
# "D1" data, so clear the blank placeholder cells the insert left behind.
$ws.Range("K2:K4").Clear()

# New column headers.
$ws.Range("K1").Value = "D1"
$ws.Range("M1").Value = "D3"

# Newly-added "D3" data (row 2 has no reading for this replicate).
$ws.Range("M3").Value = 39.889302052970102
$ws.Range("M4").Value = 38.372354057012799

# Match the saved view state (zoom + selection) from the edited workbook.
$excel.ActiveWindow.Zoom = 132
$ws.Range("L4").Select()
